# GlobalConstant.xlsx - add SpellGachaLevelAccumulatedCount global constant string
#
# Target: insert a new row into "GlobalConstantStringTable" (2nd sheet) right
# before the "OfficialCafe" row, containing the id and accumulated-count
# value, then make that sheet the active tab (it becomes the one the user is
# looking at when the file is saved).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row 4 (this shifts the former rows 4-5 down to 5-6).
$ws2.Rows.Item(4).Insert()

# Fill in the new id / value pair.
$ws2.Range("A4").Value = "SpellGachaLevelAccumulatedCount"
$ws2.Range("B4").Value = "0, 10, 45, 190, 780, 2415, 5950, 10950, 15950, 21950, 28950, 36950, 45950, 54950"

# The string table sheet becomes the active/selected sheet/tab.
$ws2.Activate()
$ws2.Range("A1").Select()
